$d = $word.ActiveDocument

# This document's body content is cyclically rotated among several paragraphs
# (the section headings themselves stay where they are). Because the moves
# form a dependency cycle (e.g. paragraph A's old text becomes paragraph B's
# new text, while paragraph B's old text is needed somewhere else), a plain
# sequence of Find/Replace calls would clobber data that a later step still
# needs. Instead we do this in two phases:
#
#   Phase 1: at every location, replace its old text with a unique
#            placeholder token that identifies that very location (e.g.
#            "%%SLOT_Docente%%" is written into the Docente paragraph).
#            A Find/Replace always mutates the paragraph the search text
#            currently lives in, so this just "tags" each paragraph.
#   Phase 2: for every location, find its own tag (it is still sitting right
#            there, since phase 1 only changed text in place) and replace it
#            with the real final text for that location (which is simply the
#            literal old text of whichever location feeds it, captured as a
#            PowerShell string literal, so it is unaffected by any mutation
#            phase 1 performed elsewhere).
#
# This guarantees correctness no matter what order Word performs the
# Find/Replace calls in.

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        throw "Find/Replace failed for: $old"
    }
}

$BR = [char]11

# The two multi-run text blocks, reconstructed with manual line breaks (and,
# for the bibliography, blank-line pairs of breaks) exactly as they appear in
# the source document.
$programaText = "1 - FUNDAMENTOS: A Engenharia e o Meio Ambiente; Os Ecossistemas. A crise energética. Fontes alternativas de energia. A sustentabilidade do meio ambiente. " + $BR + "2 - O MEIO AMBIENTE AQUÁTICO: Composição e Propriedades; Necessidade e Utilização; Requisitos de Qualidade; Poluição." + $BR + "3 - O MEIO AMBIENTE TERRESTRE: Composição e Propriedades; Necessidades e Utilização; Requisitos de Qualidade; Poluição." + $BR + "4 - O MEIO AMBIENTE ATMOSFÉRICO: Composição e Propriedades; Requisitos de Qualidade; Poluição."

$bibliografiaText = "Braga, B.P.F., M.T.,Conejo, J.G., Porto, M.F., Veras M.S., Nucci, N., Juliano, N. e Eiger, S. Introdução à Engenharia Ambiental, Makron Books, São Paulo, 1998" + $BR + $BR + "Sperling, M.V. Princípios do Tratamento Biológico de Águas Residuárias. Desa-UFMG, Minas Gerais, 1996." + $BR + $BR + "BRAGA, B.et al. Introdução à Engenharia Ambiental. São Paulo: Prentice Hall, 2002, 305 p." + $BR + $BR + "VON SPERLING, M. Introdução à qualidade das águas e ao tratamento de esgotos. 2. ed. Belo Horizonte: UFMG, 1996."

# ---------------------------------------------------------------------
# Phase 1: tag every location with a placeholder unique to that location.
# ---------------------------------------------------------------------

Replace-Text "Apresentar aos alunos os princípios fundamentais de engenharia do meio ambiente." "%%SLOT_Objetivos%%"
Replace-Text "5840671 - Francisco José Moreira Chaves" "%%SLOT_Docente%%"
Replace-Text "1 - Fundamentos da Engenharia e o Meio Ambiente. 2 - O meio ambiente aquático. 3 - O meio ambiente terrestre. 4 - O meio ambiente atmosférico ." "%%SLOT_ProgramaResumido%%"
Replace-Text $programaText "%%SLOT_Programa%%"
Replace-Text "Duas Provas  P1  1º bimestre e P2  2º bimestre" "%%SLOT_Metodo%%"
Replace-Text "MF = (P1+ P2)/2" "%%SLOT_Criterio%%"
Replace-Text "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação" "%%SLOT_Norma%%"
Replace-Text $bibliografiaText "%%SLOT_Bibliografia%%"

# ---------------------------------------------------------------------
# Phase 2: give each location its real final text.
#   Objetivos        <- old Programa resumido text
#   Docente           <- old Objetivos text
#   Programa resumido <- old Programa text
#   Programa          <- old Método value
#   Método            <- old Critério value
#   Critério          <- old Norma de recuperação value
#   Norma recuperação <- old Bibliografia text
#   Bibliografia      <- old Docente text
# ---------------------------------------------------------------------

Replace-Text "%%SLOT_Objetivos%%" "1 - Fundamentos da Engenharia e o Meio Ambiente. 2 - O meio ambiente aquático. 3 - O meio ambiente terrestre. 4 - O meio ambiente atmosférico ."
Replace-Text "%%SLOT_Docente%%" "Apresentar aos alunos os princípios fundamentais de engenharia do meio ambiente."
Replace-Text "%%SLOT_ProgramaResumido%%" $programaText
Replace-Text "%%SLOT_Programa%%" "Duas Provas  P1  1º bimestre e P2  2º bimestre"
Replace-Text "%%SLOT_Metodo%%" "MF = (P1+ P2)/2"
Replace-Text "%%SLOT_Criterio%%" "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"
Replace-Text "%%SLOT_Norma%%" $bibliografiaText
Replace-Text "%%SLOT_Bibliografia%%" "5840671 - Francisco José Moreira Chaves"
